# Commit message: "add regression for data ex china, fix Korea in WEO"
#
# The sheet lists countries ranked by an index value (col A = rank,
# col B = country, col C = avgH1, col D = gdpH1). This change:
#
#  1. Inserts "South Korea" as the new rank 28, in the row that used to
#     be occupied by "Spain" (row 29). Every country from Spain onward
#     shifts down one row AND is renumbered up by one rank (Spain
#     becomes rank 29, Sweden becomes rank 30, ..., United Kingdom
#     becomes rank 33).
#  2. Appends a brand new trailing row (rank 34) that restates
#     "United States" with a refreshed avgH1 figure (the regression
#     run excluding China).
#
# Throughout, column A values are kept as *text* (shared strings), to
# match how the rest of column A is already stored, and no new cell
# styles/number formats are introduced (mirrors the fact that
# styles.xml does not change in the target edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: insert a new blank row before the current "Spain" row (row 29).
# Spain..United States (old rows 29-34) shift down to rows 30-35, but
# keep their *old* rank numbers for now (e.g. Spain now sits in row 30
# but still shows rank "28"), so we still need to renumber them below.
# ---------------------------------------------------------------------
$ws.Rows.Item(29).Insert()

# Renumber column A for rows 29-34 (ranks 28-33) by cascading a copy of
# each row's current value (which is one rank too low, since it just
# came from the row below before the shift) from the row immediately
# below it. Doing this in ascending row order means each source cell is
# read before it gets overwritten. PasteSpecial with "values" (-4163 =
# xlPasteValues) keeps the text data type without pulling in any
# formatting/style.
for ($r = 29; $r -le 34; $r++) {
    $src = $r + 1
    $ws.Range("A$src").Copy()
    $ws.Range("A$r").PasteSpecial(-4163)
}
$ws.Application.CutCopyMode = 0

# Column B/C/D for the new South Korea row (rank 28, row 29).
$ws.Range("B29").Value = "South Korea"
$ws.Range("C29").Value = 44.2589010989011
$ws.Range("D29").Value = -4.395130442387327

# ---------------------------------------------------------------------
# Step 2: append a brand-new trailing row (row 35, rank "34") restating
# United States with an updated avgH1 figure.
# ---------------------------------------------------------------------

# Column A needs the text value "34", which doesn't exist anywhere in
# the sheet yet. Build it as text via a formula that concatenates two
# existing text cells ("3" from A4 and "4" from A5), then paste just
# the resulting value (no formula, no formatting) into A35 using a
# scratch cell. This avoids ever touching NumberFormat, which would
# otherwise register an unwanted new cell style.
$ws.Range("Z1").Formula = "=A4&A5"
$ws.Range("Z1").Copy()
$ws.Range("A35").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Application.CutCopyMode = 0

$ws.Range("B35").Value = "United States"
$ws.Range("C35").Value = 44.614688990255026
$ws.Range("D35").Value = -10.24086007423928
